# aclListSelect.xlsx — "update files to console"
#
# The diff shows two of the English shared strings in the "ch" sheet
# (column C, rows 3 & 4) were split into two runs so the trailing words
# could be highlighted in red:
#   C3: "No ACL Data"          -> "No ACL" + " data"        (red)
#   C4: "'Got ACL List Data"   -> "'Got ACL" + " list data" (red)
# (note the leading char of C4's text is U+2018 LEFT SINGLE QUOTATION MARK)
# It also moves the sheet's active selection to B19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$leftQuote = [char]0x2018

# --- C3: "No ACL Data" -> "No ACL" + red " data" -----------------------
$c3 = $ws.Range("C3")
$c3.Value = "No ACL data"
# characters 7-11 ("_data", 1-based) painted red
$c3.Characters(7, 5).Font.Color = 255

# --- C4: "'Got ACL List Data" -> "'Got ACL" + red " list data" ---------
$c4 = $ws.Range("C4")
$c4.Value = $leftQuote + "Got ACL list data"
# characters 9-18 (" list data", 1-based) painted red
$c4.Characters(9, 10).Font.Color = 255

# --- move the saved selection to B19, like the authored workbook -------
[void]$ws.Range("B19").Select()

Write-Output "aclListSelect edits applied"
